$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 348
$ws.Range("I2").Value = 348
$ws.Range("K2").Value = 348
$ws.Range("M2").Value = -235
$ws.Range("H20").Value = 5152.125
$ws.Range("I20").Value = 5152.125
$ws.Range("K20").Value = 5152.125
$ws.Range("M20").Value = -4922.125
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("M32").ClearContents()
$ws.Range("H35").Value = 5152.125
$ws.Range("I35").Value = 5152.125
$ws.Range("K35").Value = 5152.125
$ws.Range("M35").Value = -4773.125
$ws.Range("H44").Value = 43980
$ws.Range("J44").Value = 43980
$ws.Range("L44").Value = 43980
$ws.Range("N44").Value = -44904
$ws.Range("H62").Value = 8098.9165
$ws.Range("I62").Value = 6944.778
$ws.Range("J62").Value = 11561.333
$ws.Range("K62").Value = 6944.778
$ws.Range("L62").Value = 11561.333
$ws.Range("M62").Value = -6320.778
$ws.Range("N62").Value = -12809.333
$ws.Range("H65").Value = 8098.9165
$ws.Range("I65").Value = 6944.778
$ws.Range("J65").Value = 11561.333
$ws.Range("K65").Value = 34723.89
$ws.Range("L65").Value = 57806.665
$ws.Range("M65").Value = -31603.89
$ws.Range("N65").Value = -64046.665
$ws.Range("H92").Value = 18519786
$ws.Range("I92").Value = 23810800
$ws.Range("J92").Value = 1235.8334
$ws.Range("K92").Value = 23810800
$ws.Range("L92").Value = 1235.8334
$ws.Range("M92").Value = -23809552
$ws.Range("N92").Value = -3731.8334
$ws.Range("H94").Value = 3233.3333
$ws.Range("I94").Value = 3233.3333
$ws.Range("K94").Value = 3233.3333
$ws.Range("M94").Value = -2782.3333
$ws.Range("H100").Value = 3427.762
$ws.Range("J100").Value = 3700.7334
$ws.Range("L100").Value = 3700.7334
$ws.Range("N100").Value = -4782.7334
$ws.Range("H111").Value = 3034.4285
$ws.Range("I111").Value = 2987.25
$ws.Range("K111").Value = 8961.75
$ws.Range("M111").Value = -5894.75
$ws.Range("H138").Value = 3789.47
$ws.Range("I138").Value = 2658.8948
$ws.Range("J138").Value = 4054.6667
$ws.Range("K138").Value = 7976.6844
$ws.Range("L138").Value = 12164.0001
$ws.Range("M138").Value = -2836.6844
$ws.Range("N138").Value = -22444.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38465892
$ws.Range("I32").Value = 58826956
$ws.Range("J32").Value = 6111.778
$ws.Range("K32").Value = 58826956
$ws.Range("L32").Value = 6111.778
$ws.Range("M32").Value = -58826669
$ws.Range("N32").Value = -6685.778
$ws.Range("H45").Value = 2859.3171
$ws.Range("I45").Value = 2608
$ws.Range("K45").Value = 2608
$ws.Range("M45").Value = -2231
$ws.Range("H61").Value = 2651.6191
$ws.Range("I61").Value = 2645.6
$ws.Range("J61").Value = 2666.6667
$ws.Range("K61").Value = 2645.6
$ws.Range("L61").Value = 2666.6667
$ws.Range("M61").Value = -2433.6
$ws.Range("N61").Value = -3090.6667
$ws.Range("H63").Value = 5187.5
$ws.Range("I63").Value = 4875
$ws.Range("K63").Value = 4875
$ws.Range("M63").Value = -4189
$ws.Range("H66").Value = 5187.5
$ws.Range("I66").Value = 4875
$ws.Range("K66").Value = 24375
$ws.Range("M66").Value = -20943
$ws.Range("H103").Value = 122999
$ws.Range("J103").Value = 122999
$ws.Range("L103").Value = 122999
$ws.Range("N103").Value = -125343
$ws.Range("H132").Value = 3168.0312
$ws.Range("I132").Value = 2744.375
$ws.Range("J132").Value = 4439
$ws.Range("K132").Value = 8233.125
$ws.Range("L132").Value = 13317
$ws.Range("M132").Value = -5703.125
$ws.Range("N132").Value = -18377
$ws.Range("H136").Value = 2651.6191
$ws.Range("I136").Value = 2645.6
$ws.Range("J136").Value = 2666.6667
$ws.Range("K136").Value = 7936.799999999999
$ws.Range("L136").Value = 8000.000100000001
$ws.Range("M136").Value = -5386.799999999999
$ws.Range("N136").Value = -13100.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6131.6665
$ws.Range("I20").Value = 4400
$ws.Range("J20").Value = 6997.5
$ws.Range("K20").Value = 4400
$ws.Range("L20").Value = 6997.5
$ws.Range("M20").Value = -4153
$ws.Range("N20").Value = -7491.5
$ws.Range("H94").Value = 1294.6072
$ws.Range("I94").Value = 647.3158
$ws.Range("K94").Value = 647.3158
$ws.Range("M94").Value = -196.3158
$ws.Range("H99").Value = 2517.4119
$ws.Range("I99").Value = 1889.2
$ws.Range("J99").Value = 3414.8572
$ws.Range("K99").Value = 1889.2
$ws.Range("L99").Value = 3414.8572
$ws.Range("M99").Value = -391.2
$ws.Range("N99").Value = -6410.8572
$ws.Range("H105").Value = 1718.1613
$ws.Range("I105").Value = 1613.3572
$ws.Range("K105").Value = 1613.3572
$ws.Range("M105").Value = 133.6428000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.882355
$ws.Range("I7").Value = 35.125
$ws.Range("J7").Value = 99
$ws.Range("K7").Value = 35.125
$ws.Range("L7").Value = 99
$ws.Range("M7").Value = 77.875
$ws.Range("N7").Value = -325
$ws.Range("H10").Value = 1917.4286
$ws.Range("I10").Value = 1682.8
$ws.Range("J10").Value = 2504
$ws.Range("K10").Value = 1682.8
$ws.Range("L10").Value = 2504
$ws.Range("M10").Value = -1543.8
$ws.Range("N10").Value = -2782
$ws.Range("H16").Value = 1498.3334
$ws.Range("I16").Value = 1250
$ws.Range("K16").Value = 1250
$ws.Range("M16").Value = -963
$ws.Range("H29").Value = 5055
$ws.Range("I29").Value = 110
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 110
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 183
$ws.Range("N29").Value = -10586
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("K47").Value = 20000
$ws.Range("M47").Value = -19434
$ws.Range("H113").Value = 1498.3334
$ws.Range("I113").Value = 1250
$ws.Range("K113").Value = 1250
$ws.Range("M113").Value = 920
$ws.Range("H122").Value = 5888642.5
$ws.Range("I122").Value = 8338914.5
$ws.Range("K122").Value = 25016743.5
$ws.Range("M122").Value = -25014293.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.75
$ws.Range("I2").Value = 63.5
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 381
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -268
$ws.Range("N2").Value = -466
$ws.Range("H108").Value = 803.6667
$ws.Range("I108").Value = 364.4
$ws.Range("K108").Value = 1093.2
$ws.Range("M108").Value = 1786.8
$ws.Range("H121").Value = 4535786
$ws.Range("I121").Value = 707.4
$ws.Range("K121").Value = 2122.2
$ws.Range("M121").Value = -812.1999999999998
$ws.Range("H129").Value = 1766.8
$ws.Range("I129").Value = 713
$ws.Range("K129").Value = 2139
$ws.Range("M129").Value = 2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3555866.8
$ws.Range("J11").Value = 1026000
$ws.Range("L11").Value = 1026000
$ws.Range("N11").Value = -1026278
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H70").Value = 5596.5947
$ws.Range("J70").Value = 4909.091
$ws.Range("L70").Value = 4909.091
$ws.Range("N70").Value = -5449.091
$ws.Range("H73").Value = 5596.5947
$ws.Range("J73").Value = 4909.091
$ws.Range("L73").Value = 4909.091
$ws.Range("N73").Value = -6781.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9624.833000000001
$ws.Range("I122").Value = 5437.5
$ws.Range("K122").Value = 16312.5
$ws.Range("M122").Value = -13862.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 210000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H62").Value = 5312.6665
$ws.Range("I62").Value = 3240.3333
$ws.Range("J62").Value = 6348.8335
$ws.Range("K62").Value = 3240.3333
$ws.Range("L62").Value = 6348.8335
$ws.Range("M62").Value = -2616.3333
$ws.Range("N62").Value = -7596.8335
$ws.Range("H65").Value = 5312.6665
$ws.Range("I65").Value = 3240.3333
$ws.Range("J65").Value = 6348.8335
$ws.Range("K65").Value = 16201.6665
$ws.Range("L65").Value = 31744.1675
$ws.Range("M65").Value = -13081.6665
$ws.Range("N65").Value = -37984.1675
$ws.Range("H122").Value = 43488212
$ws.Range("I122").Value = 50008908
$ws.Range("J122").Value = 16896.334
$ws.Range("K122").Value = 150026724
$ws.Range("L122").Value = 50689.00199999999
$ws.Range("M122").Value = -150024274
$ws.Range("N122").Value = -55589.00199999999
$ws.Range("H132").Value = 2286.3696
$ws.Range("I132").Value = 1993.6666
$ws.Range("J132").Value = 5359.75
$ws.Range("K132").Value = 5980.9998
$ws.Range("L132").Value = 16079.25
$ws.Range("M132").Value = -3450.9998
$ws.Range("N132").Value = -21139.25
